{"js": "// Insert a new \"-could maybe add a healing block...\" paragraph right\n// after the existing \"-could maybe add a timer to D.Vas menu...\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText =\n  \"-could maybe add a timer to D.Vas menu so she can\\u2019t just stay in it forever.\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate target paragraph: \" + targetText);\n}\n\nconst newText =\n  \"-could maybe add a healing block when a player is about to explode from genji\\u2019s ultimate\";\n\n// Inserting \"After\" the found paragraph clones its paragraph/run\n// formatting (the fi-FI language mark), matching how Word splits a\n// paragraph when the author pressed Enter at the end of the line.\ntarget.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"-could maybe add a healing block...\" paragraph right\n# after the existing \"-could maybe add a timer to D.Vas menu...\" paragraph.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"-could maybe add a timer to D.Vas menu so she can\" + [char]0x2019 + \"t just stay in it forever.\")\nif (-not $found) {\n    throw \"Could not locate target paragraph\"\n}\n\n# Expand the found (collapsed) range back out to the whole paragraph so\n# InsertParagraphAfter lands after the paragraph mark, not mid-sentence.\n$range.Expand(4)  # wdParagraph\n\n$range.InsertParagraphAfter()\n\n# The freshly inserted (empty) paragraph immediately follows $range.\n$newPara = $range.Next(4, 1)  # wdParagraph\n$newPara.InsertBefore(\"-could maybe add a healing block when a player is about to explode from genji\" + [char]0x2019 + \"s ultimate\")\n"}
